$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Row 599
$ws.Range("A599").Value = "Demo"
$ws.Range("B599").Value = "9/29/2016"
$ws.Range("C599").Value = "1700"
$ws.Range("D599").Value = "OSG"
$ws.Range("E599").Value = "2009"

# Row 600
$ws.Range("A600").Value = "Demo"
$ws.Range("B600").Value = "9/29/2016"
$ws.Range("C600").Value = "1730"
$ws.Range("D600").Value = "SSB"
$ws.Range("E600").Value = "E118"
$ws.Range("F600").Value = "Check with prof re :wireless mouse"

# Row 601
$ws.Range("A601").Value = "Demo"
$ws.Range("B601").Value = "9/29/2016"
$ws.Range("C601").Value = "1730"
$ws.Range("D601").Value = "SSB"
$ws.Range("E601").Value = "E118"
$ws.Range("F601").Value = "Check with prof re :wireless mouse"

# Row 602
$ws.Range("A602").Value = "Demo"
$ws.Range("B602").Value = "9/29/2016"
$ws.Range("C602").Value = "1600"
$ws.Range("D602").Value = "HNE"
$ws.Range("E602").Value = "401"

# Row 603
$ws.Range("A603").Value = "Demo"
$ws.Range("B603").Value = "9/29/2016"
$ws.Range("C603").Value = "1900"
$ws.Range("D603").Value = "SSB"
$ws.Range("E603").Value = "N108"
$ws.Range("F603").Value = "Client using neck mic"

# Row 604
$ws.Range("A604").Value = "Setup Mic"
$ws.Range("B604").Value = "9/29/2016"
$ws.Range("C604").Value = "1800"
$ws.Range("D604").Value = "DB"
$ws.Range("E604").Value = "2027"
$ws.Range("F604").Value = "Neck mic and small PA from DB 0003"

# Row 605
$ws.Range("A605").Value = "Pickup Mic"
$ws.Range("B605").Value = "9/29/2016"
$ws.Range("C605").Value = "2100"
$ws.Range("D605").Value = "DB"
$ws.Range("E605").Value = "2027"
$ws.Range("F605").Value = "Return neck mic and small PA to DB 0003"

# Row 606
$ws.Range("A606").Value = "Demo"
$ws.Range("B606").Value = "9/29/2016"
$ws.Range("C606").Value = "1700"
$ws.Range("D606").Value = "SSB"
$ws.Range("E606").Value = "W141"
$ws.Range("F606").Value = "PC, neck mic and podium mic"

# Row 607
$ws.Range("A607").Value = "AV Shutdown"
$ws.Range("B607").Value = "9/29/2016"
$ws.Range("C607").Value = "1930"
$ws.Range("D607").Value = "SSB"
$ws.Range("E607").Value = "W141"

# Row 611 (rows 608-610 intentionally left blank)
$ws.Range("A611").Value = "Demo"
$ws.Range("B611").Value = "10/3/2016"
$ws.Range("C611").Value = "1820"
$ws.Range("D611").Value = "DB"
$ws.Range("E611").Value = "0009"

# Row 612
$ws.Range("A612").Value = "Demo"
$ws.Range("B612").Value = "10/3/2016"
$ws.Range("C612").Value = "1900"
$ws.Range("D612").Value = "SSB"
$ws.Range("E612").Value = "S124"

# Row 613
$ws.Range("A613").Value = "Pickup Mic"
$ws.Range("B613").Value = "10/3/2016"
$ws.Range("C613").Value = "1730"
$ws.Range("D613").Value = "HNE"
$ws.Range("E613").Value = "281"
$ws.Range("F613").Value = "Return 4 IR mics, receivers, cables and stands to HNES 003"

# Match the final selection/active cell from the edit
$ws.Range("F619").Select()
